$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Indent of the existing Relative-Frequency column (C2:C9) grows from 3 to 4
#    (mutates the shared style used by those cells -> becomes cellXfs[1] w/ indent=4)
$ws.Range("C2:C9").IndentLevel = 4

# 2) New column width for column C (feature column), target raw width ~16.94
#    (COM ColumnWidth snaps to the nearest 1/6-character pixel grid, so this
#    is the closest achievable value to the authored 16.94 width)
$ws.Columns.Item(3).ColumnWidth = 16.17

# 3) Existing "Number of purchases" distribution gains one more bucket (">9")
$ws.Range("G8").Value = 106
$ws.Range("G9").Value = 80

$ws.Range("A10").Value = ">9"
$ws.Range("B10").Value = 298
$ws.Range("C10").Value = 0.003
$ws.Range("C10").HorizontalAlignment = -4108   # xlCenter

$ws.Range("G10").Value = 37
$ws.Range("G11").Value = 21
$ws.Range("G12").Value = 14
$ws.Range("G13").Value = 13
$ws.Range("G14").Value = 5
$ws.Range("G15").Value = 5
$ws.Range("G16").Value = 4
$ws.Range("G17").Value = 2
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("G21").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("G25").Value = 1

# 4) k-means feature block: wrapped placeholder cells in columns G and J
$ws.Range("G10:G25").WrapText = $true
$ws.Range("J13:J30").WrapText = $true

# 5) Selection follows the extended Relative-Frequency column
[void]$ws.Range("C4:C10").Select()
